# WB_Test_Report_2025-12-22.xlsx — batch interrupted after milestone 3 (Quote
# Rated Successfully never completed cleanly -> Test Execution Failed) so only
# suite #16 (BOP (OH)) got as far as it did before the grace-wait fallback
# kicked in. Reflect that in the report: rename the BOP sheet, update its
# single result row, and trim the milestone timeline down to what actually
# ran.

$wb = $excel.ActiveWorkbook

# --- Sheet 2 ("BOP_1" -> "BOP_16") ---------------------------------------
$bopSheet = $wb.Worksheets.Item("BOP_1")
$bopSheet.Name = "BOP_16"

# --- Summary sheet: update the single iteration row ----------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("A2").Value = 16
$summary.Range("B2").Value = "BOP (OH)"
$summary.Range("C2").Value = "N/A"
$summary.Range("F2").Value = "'228.42"

# --- BOP_16 sheet: update milestone timings, mark suite failed, and drop
#     the two milestones that never ran because the batch was interrupted --
$bopSheet.Range("C2").Value = "'84.67s"
$bopSheet.Range("D2").Value = "'2025-12-22T14:51:40.324Z"

$bopSheet.Range("C3").Value = "'138.63s"
$bopSheet.Range("D3").Value = "'2025-12-22T14:53:58.960Z"

$bopSheet.Range("A4").Value = "Test Execution Failed"
$bopSheet.Range("B4").Value = "FAILED"
$bopSheet.Range("C4").Value = "'5.12s"
$bopSheet.Range("D4").Value = "'2025-12-22T14:54:04.077Z"

# Rows 5 and 6 ("Submitting for Approval" / "Test Execution Failed") belong
# to milestones past the interruption point -- remove them so the sheet ends
# at row 4.
$bopSheet.Rows.Item(6).Delete()
$bopSheet.Rows.Item(5).Delete()
